# Auto-generated Excel COM-interop script to apply market-price data updates
# across all 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 118.6
$ws.Range("I5").Value = 98.25
$ws.Range("K5").Value = 98.25
$ws.Range("M5").Value = 16.75
$ws.Range("H28").Value = 1636.7142
$ws.Range("I28").Value = 1636.7142
$ws.Range("K28").Value = 1636.7142
$ws.Range("M28").Value = -1151.7142
$ws.Range("H62").Value = 7926.2
$ws.Range("I62").Value = 7873
$ws.Range("K62").Value = 7873
$ws.Range("M62").Value = -7249
$ws.Range("H65").Value = 7926.2
$ws.Range("I65").Value = 7873
$ws.Range("K65").Value = 39365
$ws.Range("M65").Value = -36245
$ws.Range("H98").Value = 3296.3333
$ws.Range("I98").Value = 1181.6666
$ws.Range("J98").Value = 3825
$ws.Range("K98").Value = 1181.6666
$ws.Range("L98").Value = 3825
$ws.Range("M98").Value = 316.3334
$ws.Range("N98").Value = -6821
$ws.Range("H122").Value = 3296.3333
$ws.Range("I122").Value = 1181.6666
$ws.Range("J122").Value = 3825
$ws.Range("K122").Value = 3544.9998
$ws.Range("L122").Value = 11475
$ws.Range("M122").Value = -1094.9998
$ws.Range("N122").Value = -16375

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3618.625
$ws.Range("I132").Value = 3389.8
$ws.Range("K132").Value = 10169.4
$ws.Range("M132").Value = -7639.400000000001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7047
$ws.Range("I86").Value = 600
$ws.Range("J86").Value = 10270.5
$ws.Range("K86").Value = 600
$ws.Range("L86").Value = 10270.5
$ws.Range("M86").Value = 523
$ws.Range("N86").Value = -12516.5
$ws.Range("H89").Value = 7047
$ws.Range("I89").Value = 600
$ws.Range("J89").Value = 10270.5
$ws.Range("K89").Value = 3000
$ws.Range("L89").Value = 51352.5
$ws.Range("M89").Value = 2616
$ws.Range("N89").Value = -62584.5
$ws.Range("H134").Value = 1317.091
$ws.Range("I134").Value = 1248.8
$ws.Range("K134").Value = 3746.4
$ws.Range("M134").Value = -1211.4
$ws.Range("H135").Value = 99994.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 99994.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 99994.5
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -110134.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 500
$ws.Range("J14").Value = 500
$ws.Range("L14").Value = 500
$ws.Range("N14").Value = -840
$ws.Range("H54").Value = 34506.5
$ws.Range("I54").Value = 28934.5
$ws.Range("K54").Value = 28934.5
$ws.Range("M54").Value = -28276.5
$ws.Range("H58").Value = 2250.611
$ws.Range("I58").Value = 2235.9412
$ws.Range("K58").Value = 2235.9412
$ws.Range("M58").Value = -2032.9412
$ws.Range("H92").Value = 43063
$ws.Range("J92").Value = 43063
$ws.Range("L92").Value = 43063
$ws.Range("N92").Value = -48055
$ws.Range("H95").Value = 17460
$ws.Range("J95").Value = 17460
$ws.Range("L95").Value = 17460
$ws.Range("N95").Value = -22952
$ws.Range("H96").Value = 10400
$ws.Range("J96").Value = 10400
$ws.Range("L96").Value = 10400
$ws.Range("N96").Value = -15892
$ws.Range("H105").Value = 3775
$ws.Range("I105").Value = 3079.4
$ws.Range("J105").Value = 4310.077
$ws.Range("K105").Value = 3079.4
$ws.Range("L105").Value = 4310.077
$ws.Range("M105").Value = -1332.4
$ws.Range("N105").Value = -7804.077
$ws.Range("H106").Value = 187947.33
$ws.Range("J106").Value = 187947.33
$ws.Range("L106").Value = 187947.33
$ws.Range("N106").Value = -190471.33
$ws.Range("H107").Value = 620.5333000000001
$ws.Range("I107").Value = 319
$ws.Range("J107").Value = 1072.8334
$ws.Range("K107").Value = 319
$ws.Range("L107").Value = 1072.8334
$ws.Range("M107").Value = 1601
$ws.Range("N107").Value = -4912.8334
$ws.Range("H110").Value = 75000
$ws.Range("J110").Value = 75000
$ws.Range("L110").Value = 75000
$ws.Range("N110").Value = -83180
$ws.Range("H111").Value = 47499.5
$ws.Range("J111").Value = 47499.5
$ws.Range("L111").Value = 47499.5
$ws.Range("N111").Value = -55679.5
$ws.Range("H134").Value = 2739.3333
$ws.Range("I134").Value = 2739.3333
$ws.Range("K134").Value = 8217.999899999999
$ws.Range("M134").Value = -5682.999899999999
$ws.Range("H136").Value = 2250.611
$ws.Range("I136").Value = 2235.9412
$ws.Range("K136").Value = 6707.823600000001
$ws.Range("M136").Value = -4157.823600000001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 100
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H131").Value = 2871
$ws.Range("J131").Value = 2871
$ws.Range("L131").Value = 8613
$ws.Range("N131").Value = -18693
$ws.Range("H137").Value = 2249.3572
$ws.Range("J137").Value = 3399
$ws.Range("L137").Value = 10197
$ws.Range("N137").Value = -20397

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1211.6666
$ws.Range("I113").Value = 1211.6666
$ws.Range("K113").Value = 1211.6666
$ws.Range("M113").Value = 958.3334
$ws.Range("H135").Value = 45000
$ws.Range("J135").Value = 45000
$ws.Range("L135").Value = 45000
$ws.Range("N135").Value = -55140

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1058.3
$ws.Range("J16").Value = 1322.5
$ws.Range("L16").Value = 1322.5
$ws.Range("N16").Value = -1662.5
$ws.Range("H40").Value = 4054.2222
$ws.Range("J40").Value = 4555.4287
$ws.Range("L40").Value = 4555.4287
$ws.Range("N40").Value = -4827.4287
$ws.Range("H46").Value = 40736.23
$ws.Range("I46").Value = 85428.836
$ws.Range("J46").Value = 2428.2856
$ws.Range("K46").Value = 85428.836
$ws.Range("L46").Value = 2428.2856
$ws.Range("M46").Value = -85240.836
$ws.Range("N46").Value = -2804.2856
$ws.Range("H61").Value = 3153.1
$ws.Range("I61").Value = 3059.111
$ws.Range("K61").Value = 3059.111
$ws.Range("M61").Value = -2857.111
$ws.Range("H68").Value = 1995.8334
$ws.Range("I68").Value = 1627
$ws.Range("J68").Value = 2259.2856
$ws.Range("K68").Value = 1627
$ws.Range("L68").Value = 2259.2856
$ws.Range("M68").Value = -878
$ws.Range("N68").Value = -3757.2856
$ws.Range("H71").Value = 1995.8334
$ws.Range("I71").Value = 1627
$ws.Range("J71").Value = 2259.2856
$ws.Range("K71").Value = 8135
$ws.Range("L71").Value = 11296.428
$ws.Range("M71").Value = -4391
$ws.Range("N71").Value = -18784.428
$ws.Range("H113").Value = 3153.1
$ws.Range("I113").Value = 3059.111
$ws.Range("K113").Value = 3059.111
$ws.Range("M113").Value = -889.1109999999999
$ws.Range("H132").Value = 3349.7
$ws.Range("I132").Value = 2748
$ws.Range("K132").Value = 8244
$ws.Range("M132").Value = -5714
$ws.Range("H136").Value = 3068.4119
$ws.Range("J136").Value = 4169
$ws.Range("L136").Value = 12507
$ws.Range("N136").Value = -17607

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H122").Value = 11125.5
$ws.Range("J122").Value = 11125.5
$ws.Range("L122").Value = 33376.5
$ws.Range("N122").Value = -38276.5
$ws.Range("H136").Value = 2037.7693
$ws.Range("I136").Value = 1790.9166
$ws.Range("K136").Value = 5372.7498
$ws.Range("M136").Value = -2822.7498
